# Auto-generated Excel COM-interop script to apply the Kujata_Profits.xlsx value updates.
# Updates per-row, per-column numeric values (and removes a few now-empty cells) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1914.3125
$ws.Range("J112").Value = 2655.5
$ws.Range("L112").Value = 7966.5
$ws.Range("N112").Value = -10182.5

$ws.Range("H129").Value = 841.9818
$ws.Range("I129").Value = 615.2
$ws.Range("J129").Value = 864.66
$ws.Range("K129").Value = 1845.6
$ws.Range("L129").Value = 2593.98
$ws.Range("M129").Value = 3154.4
$ws.Range("N129").Value = -12593.98

$ws.Range("H131").Value = 767
$ws.Range("I131").Value = 770.7778
$ws.Range("J131").Value = 750
$ws.Range("K131").Value = 2312.3334
$ws.Range("L131").Value = 2250
$ws.Range("M131").Value = 2727.6666
$ws.Range("N131").Value = -12330

$ws.Range("H132").Value = 8341583
$ws.Range("I132").Value = 9529132
$ws.Range("J132").Value = 28741.2
$ws.Range("K132").Value = 28587396
$ws.Range("L132").Value = 86223.6
$ws.Range("M132").Value = -28584866
$ws.Range("N132").Value = -91283.6

$ws.Range("H133").Value = 37375
$ws.Range("J133").Value = 37375
$ws.Range("L133").Value = 37375
$ws.Range("N133").Value = -47495

$ws.Range("H134").Value = 37232
$ws.Range("J134").Value = 37232
$ws.Range("L134").Value = 37232
$ws.Range("N134").Value = -47372

$ws.Range("H136").Value = 36843
$ws.Range("J136").Value = 36843
$ws.Range("L136").Value = 36843
$ws.Range("N136").Value = -47043

$ws.Range("H137").Value = 2170.9607
$ws.Range("I137").Value = 1822.3334
$ws.Range("J137").Value = 2480.8518
$ws.Range("K137").Value = 5467.0002
$ws.Range("L137").Value = 7442.555399999999
$ws.Range("M137").Value = -2917.0002
$ws.Range("N137").Value = -12542.5554

$ws.Range("H139").Value = 94340
$ws.Range("J139").Value = 94340
$ws.Range("L139").Value = 94340
$ws.Range("N139").Value = -104620

$ws.Range("H140").Value = 37426
$ws.Range("J140").Value = 37426
$ws.Range("L140").Value = 37426
$ws.Range("N140").Value = -47786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10401.393
$ws.Range("I32").Value = 7722.7354
$ws.Range("K32").Value = 7722.7354
$ws.Range("M32").Value = -7435.7354

$ws.Range("H97").Value = 559
$ws.Range("I97").Value = 541.13336
$ws.Range("J97").Value = 693
$ws.Range("K97").Value = 541.13336
$ws.Range("L97").Value = 693
$ws.Range("M97").Value = -45.13336000000004
$ws.Range("N97").Value = -1685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4159.9375
$ws.Range("I86").Value = 4684.778
$ws.Range("J86").Value = 3485.1428
$ws.Range("K86").Value = 4684.778
$ws.Range("L86").Value = 3485.1428
$ws.Range("M86").Value = -3561.778
$ws.Range("N86").Value = -5731.1428

$ws.Range("H89").Value = 4159.9375
$ws.Range("I89").Value = 4684.778
$ws.Range("J89").Value = 3485.1428
$ws.Range("K89").Value = 23423.89
$ws.Range("L89").Value = 17425.714
$ws.Range("M89").Value = -17807.89
$ws.Range("N89").Value = -28657.714

$ws.Range("H134").Value = 3634.4595
$ws.Range("I134").Value = 884.8214
$ws.Range("J134").Value = 12188.889
$ws.Range("K134").Value = 2654.4642
$ws.Range("L134").Value = 36566.667
$ws.Range("M134").Value = -119.4642000000003
$ws.Range("N134").Value = -41636.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1523.0193
$ws.Range("I31").Value = 1419.7142
$ws.Range("J31").Value = 1956.9
$ws.Range("K31").Value = 1419.7142
$ws.Range("L31").Value = 1956.9
$ws.Range("M31").Value = -1124.7142
$ws.Range("N31").Value = -2546.9

$ws.Range("H34").Value = 1523.0193
$ws.Range("I34").Value = 1419.7142
$ws.Range("J34").Value = 1956.9
$ws.Range("K34").Value = 1419.7142
$ws.Range("L34").Value = 1956.9
$ws.Range("M34").Value = -1217.7142
$ws.Range("N34").Value = -2360.9

$ws.Range("H132").Value = 2254.2942
$ws.Range("I132").Value = 2407.0833
$ws.Range("K132").Value = 7221.249899999999
$ws.Range("M132").Value = -4691.249899999999

$ws.Range("H141").Value = 1024296
$ws.Range("J141").Value = 1024296
$ws.Range("L141").Value = 1024296
$ws.Range("N141").Value = -1034656

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 397.7143
$ws.Range("I33").Value = 263.33334
$ws.Range("J33").Value = 498.5
$ws.Range("K33").Value = 1580.00004
$ws.Range("L33").Value = 2991
$ws.Range("M33").Value = -1297.00004
$ws.Range("N33").Value = -3557

$ws.Range("H38").Value = 360.66666

$ws.Range("H68").Value = 681.8125
$ws.Range("I68").Value = 797.8571
$ws.Range("J68").Value = 591.55554
$ws.Range("K68").Value = 2393.5713
$ws.Range("L68").Value = 1774.66662
$ws.Range("M68").Value = -1582.5713
$ws.Range("N68").Value = -3396.66662

$ws.Range("H71").Value = 681.8125
$ws.Range("I71").Value = 797.8571
$ws.Range("J71").Value = 591.55554
$ws.Range("K71").Value = 7180.7139
$ws.Range("L71").Value = 5323.99986
$ws.Range("M71").Value = -3124.7139
$ws.Range("N71").Value = -13435.99986

$ws.Range("H82").Value = 10180
$ws.Range("I82").Value = 2900
$ws.Range("K82").Value = 8700
$ws.Range("M82").Value = -8294

$ws.Range("H85").Value = 10180
$ws.Range("I85").Value = 2900
$ws.Range("K85").Value = 8700
$ws.Range("M85").Value = -7296

$ws.Range("H107").Value = 7691.2856
$ws.Range("J107").Value = 9636.182
$ws.Range("L107").Value = 28908.546
$ws.Range("N107").Value = -32748.546

$ws.Range("H113").Value = 729.02856
$ws.Range("J113").Value = 746.3929
$ws.Range("L113").Value = 2239.1787
$ws.Range("N113").Value = -6579.1787

$ws.Range("H122").Value = 930.85364
$ws.Range("J122").Value = 1066.4482
$ws.Range("L122").Value = 9598.033800000001
$ws.Range("N122").Value = -14498.0338

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5753631.5
$ws.Range("I12").Value = 6018176.5
$ws.Range("J12").Value = 3505000
$ws.Range("K12").Value = 6018176.5
$ws.Range("L12").Value = 3505000
$ws.Range("M12").Value = -6018036.5
$ws.Range("N12").Value = -3505280

$ws.Range("H122").Value = 1100.4286
$ws.Range("I122").Value = 1100.4286
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3301.2858
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -851.2857999999997
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 2194.5454
$ws.Range("I126").Value = 1917.5
$ws.Range("K126").Value = 5752.5
$ws.Range("M126").Value = -3282.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 2253
$ws.Range("I13").Value = 2253
$ws.Range("K13").Value = 2253
$ws.Range("M13").Value = -2113

$ws.Range("H40").Value = 2326.1904
$ws.Range("I40").Value = 1422.6666
$ws.Range("J40").Value = 4585
$ws.Range("K40").Value = 1422.6666
$ws.Range("L40").Value = 4585
$ws.Range("M40").Value = -1286.6666
$ws.Range("N40").Value = -4857

$ws.Range("H68").Value = 1301.875
$ws.Range("I68").Value = 1302
$ws.Range("K68").Value = 1302
$ws.Range("M68").Value = -553

$ws.Range("H71").Value = 1301.875
$ws.Range("I71").Value = 1302
$ws.Range("K71").Value = 6510
$ws.Range("M71").Value = -2766

$ws.Range("H93").Value = 1032.75
$ws.Range("I93").Value = 1032.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1032.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 215.25
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100002300
$ws.Range("I62").Value = 250001250
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 250001250
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -250000626
$ws.Range("N62").Value = -4248

$ws.Range("H65").Value = 100002300
$ws.Range("I65").Value = 250001250
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 1250006250
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -1250003130
$ws.Range("N65").Value = -21240

$ws.Range("H107").Value = 565.8333
$ws.Range("J107").Value = 699.5
$ws.Range("L107").Value = 2098.5
$ws.Range("N107").Value = -5938.5

